$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2, 4) '34.785.95'
Set-TextValue $ws.Cells.Item(2, 5) '  -0.23%  '
Set-TextValue $ws.Cells.Item(3, 4) '1.822.47'
Set-TextValue $ws.Cells.Item(3, 5) '  +0.71%  '
Set-TextValue $ws.Cells.Item(4, 5) '  +0.53%  '
Set-TextValue $ws.Cells.Item(5, 4) '230.77'
Set-TextValue $ws.Cells.Item(5, 5) '  -0.47%  '
Set-TextValue $ws.Cells.Item(6, 4) '0.614'
Set-TextValue $ws.Cells.Item(6, 5) '  +1.13%  '
Set-TextValue $ws.Cells.Item(7, 5) '  +0.51%  '
Set-TextValue $ws.Cells.Item(8, 4) '39.64'
Set-TextValue $ws.Cells.Item(8, 5) '  -1.26%  '
Set-TextValue $ws.Cells.Item(9, 4) '0.319'
Set-TextValue $ws.Cells.Item(9, 5) '  +2.08%  '
Set-TextValue $ws.Cells.Item(10, 5) '  -0.18%  '
Set-TextValue $ws.Cells.Item(11, 4) '0.0990'
Set-TextValue $ws.Cells.Item(11, 5) '  -0.53%  '
Set-TextValue $ws.Cells.Item(12, 4) '2.086.53'
Set-TextValue $ws.Cells.Item(12, 5) '  +0.75%  '
Set-TextValue $ws.Cells.Item(13, 2) 'WrappedEther'
Set-TextValue $ws.Cells.Item(13, 3) 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Cells.Item(13, 4) '1.836.80'
Set-TextValue $ws.Cells.Item(13, 5) '  +1.44%  '
Set-TextValue $ws.Cells.Item(14, 2) 'Chainlink'
Set-TextValue $ws.Cells.Item(14, 3) 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Cells.Item(14, 4) '11.26'
Set-TextValue $ws.Cells.Item(14, 5) '  +1.85%  '
Set-TextValue $ws.Cells.Item(15, 4) '0.665'
Set-TextValue $ws.Cells.Item(15, 5) '  +1.48%  '
Set-TextValue $ws.Cells.Item(16, 4) '4.62'
Set-TextValue $ws.Cells.Item(16, 5) '  -1.08%  '
Set-TextValue $ws.Cells.Item(17, 4) '34.652.83'
Set-TextValue $ws.Cells.Item(17, 5) '  -0.50%  '
Set-TextValue $ws.Cells.Item(18, 4) '69.47'
Set-TextValue $ws.Cells.Item(18, 5) '  +0.75%  '
Set-TextValue $ws.Cells.Item(19, 4) '0.0₃0784'
Set-TextValue $ws.Cells.Item(19, 5) '  +0.19%  '
Set-TextValue $ws.Cells.Item(20, 4) '239.14'
Set-TextValue $ws.Cells.Item(20, 5) '  +1.01%  '
Set-TextValue $ws.Cells.Item(21, 4) '12.04'
Set-TextValue $ws.Cells.Item(21, 5) '  +2.47%  '
Set-TextValue $ws.Cells.Item(22, 5) '  +0.12%  '
Set-TextValue $ws.Cells.Item(23, 5) '  +0.42%  '
Set-TextValue $ws.Cells.Item(24, 5) '  -0.14%  '
Set-TextValue $ws.Cells.Item(25, 4) '173.22'
Set-TextValue $ws.Cells.Item(25, 5) '  +0.26%  '
Set-TextValue $ws.Cells.Item(26, 5) '  -2.65%  '
Set-TextValue $ws.Cells.Item(27, 4) '0.124'
Set-TextValue $ws.Cells.Item(27, 5) '  +2.92%  '
Set-TextValue $ws.Cells.Item(28, 4) '17.27'
Set-TextValue $ws.Cells.Item(28, 5) '  -0.64%  '
Set-TextValue $ws.Cells.Item(29, 5) '  -6.87%  '
Set-TextValue $ws.Cells.Item(30, 5) '  +0.67%  '
Set-TextValue $ws.Cells.Item(31, 4) '0.0547'
Set-TextValue $ws.Cells.Item(31, 5) '  -0.41%  '
Set-TextValue $ws.Cells.Item(32, 4) '3.90'
Set-TextValue $ws.Cells.Item(32, 5) '  +0.80%  '
Set-TextValue $ws.Cells.Item(33, 4) '3.91'
Set-TextValue $ws.Cells.Item(33, 5) '  -0.95%  '
Set-TextValue $ws.Cells.Item(34, 4) '1.22'
Set-TextValue $ws.Cells.Item(34, 5) '  +5.07%  '
Set-TextValue $ws.Cells.Item(35, 5) '  +1.71%  '
Set-TextValue $ws.Cells.Item(36, 5) '  +12.24%  '
Set-TextValue $ws.Cells.Item(37, 4) '0.694'
Set-TextValue $ws.Cells.Item(37, 5) '  +2.83%  '
Set-TextValue $ws.Cells.Item(38, 4) '91.34'
Set-TextValue $ws.Cells.Item(38, 5) '  -1.78%  '
Set-TextValue $ws.Cells.Item(39, 4) '1.338.60'
Set-TextValue $ws.Cells.Item(39, 5) '  +2.69%  '
Set-TextValue $ws.Cells.Item(40, 5) '  +2.67%  '
Set-TextValue $ws.Cells.Item(41, 4) '0.0192'
Set-TextValue $ws.Cells.Item(41, 5) '  +0.27%  '
Set-TextValue $ws.Cells.Item(42, 4) '14.42'
Set-TextValue $ws.Cells.Item(42, 5) '  -2.23%  '
Set-TextValue $ws.Cells.Item(43, 4) '2.43'
Set-TextValue $ws.Cells.Item(43, 5) '  -0.30%  '
Set-TextValue $ws.Cells.Item(44, 5) '  -3.84%  '
Set-TextValue $ws.Cells.Item(45, 4) '2.74'
Set-TextValue $ws.Cells.Item(45, 5) '  -0.23%  '
Set-TextValue $ws.Cells.Item(46, 4) '6.25'
Set-TextValue $ws.Cells.Item(46, 5) '  +0.03%  '
Set-TextValue $ws.Cells.Item(47, 4) '0.0522'
Set-TextValue $ws.Cells.Item(47, 5) '  +1.97%  '
Set-TextValue $ws.Cells.Item(48, 4) '2.002.77'
Set-TextValue $ws.Cells.Item(48, 5) '  +0.89%  '
Set-TextValue $ws.Cells.Item(49, 5) '  +0.62%  '
Set-TextValue $ws.Cells.Item(50, 4) '0.0668'
Set-TextValue $ws.Cells.Item(50, 5) '  +3.82%  '
Set-TextValue $ws.Cells.Item(51, 4) '98.03'
Set-TextValue $ws.Cells.Item(51, 5) '  -1.01%  '
